$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new blank column before the
#     old "Late" column (was column N), shifting Late / heading / Outstanding
#     one column to the right (N->O, O->P, P->Q). ---
$ws = $wb.Worksheets.Item("Repayment schedule")

# Capture the width of column M (last untouched numeric column) so the
# newly inserted column N can inherit the same display width.
$mWidth = $ws.Columns.Item(13).ColumnWidth

# Insert a new column at N; everything from N onward shifts right.
$ws.Columns.Item(14).Insert()

# Match the new column's width to its left neighbour (column M), matching
# the "11" character width Excel used when the column was created.
$ws.Columns.Item(14).ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab (this also clears the
# tabSelected flag previously held by "Edit Repayment Schedule") and move
# the selection cursor to where the author left it.
$ws.Activate() | Out-Null
$ws.Range("S7").Select() | Out-Null
